$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 603.4545000000001
$ws.Range("I19").Value = 641.7143
$ws.Range("J19").Value = 536.5
$ws.Range("K19").Value = 641.7143
$ws.Range("L19").Value = 536.5
$ws.Range("M19").Value = -466.7143
$ws.Range("N19").Value = -886.5

$ws.Range("H28").Value = 766.2174
$ws.Range("I28").Value = 779.9474
$ws.Range("J28").Value = 701
$ws.Range("K28").Value = 779.9474
$ws.Range("L28").Value = 701
$ws.Range("M28").Value = -294.9474
$ws.Range("N28").Value = -1671

$ws.Range("H62").Value = 50004640
$ws.Range("I62").Value = 62504364
$ws.Range("J62").Value = 5753
$ws.Range("K62").Value = 62504364
$ws.Range("L62").Value = 5753
$ws.Range("M62").Value = -62503740
$ws.Range("N62").Value = -7001

$ws.Range("H65").Value = 50004640
$ws.Range("I65").Value = 62504364
$ws.Range("J65").Value = 5753
$ws.Range("K65").Value = 312521820
$ws.Range("L65").Value = 28765
$ws.Range("M65").Value = -312518700
$ws.Range("N65").Value = -35005

$ws.Range("H107").Value = 2326
$ws.Range("I107").Value = 2057.257
$ws.Range("J107").Value = 3266.6
$ws.Range("K107").Value = 2057.257
$ws.Range("L107").Value = 3266.6
$ws.Range("M107").Value = -137.2570000000001
$ws.Range("N107").Value = -7106.6

$ws.Range("H116").Value = 6933.212
$ws.Range("I116").Value = 6084.3887
$ws.Range("J116").Value = 7951.8
$ws.Range("K116").Value = 6084.3887
$ws.Range("L116").Value = 7951.8
$ws.Range("M116").Value = -2642.3887
$ws.Range("N116").Value = -14835.8

$ws.Range("H120").Value = 69761
$ws.Range("J120").Value = 69761
$ws.Range("L120").Value = 69761
$ws.Range("N120").Value = -79437

$ws.Range("H132").Value = 2593.9678
$ws.Range("I132").Value = 2575.7778
$ws.Range("K132").Value = 7727.3334
$ws.Range("M132").Value = -5197.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6098.4644
$ws.Range("I32").Value = 6301.8335
$ws.Range("J32").Value = 4878.25
$ws.Range("K32").Value = 6301.8335
$ws.Range("L32").Value = 4878.25
$ws.Range("M32").Value = -6014.8335
$ws.Range("N32").Value = -5452.25

$ws.Range("H110").Value = 2342.6333
$ws.Range("I110").Value = 1906.8695
$ws.Range("K110").Value = 1906.8695
$ws.Range("M110").Value = 138.1305

$ws.Range("H132").Value = 2800.255
$ws.Range("I132").Value = 2724.5527
$ws.Range("K132").Value = 8173.658100000001
$ws.Range("M132").Value = -5643.658100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 333373340
$ws.Range("I86").Value = 333373340
$ws.Range("K86").Value = 333373340
$ws.Range("M86").Value = -333372217

$ws.Range("H89").Value = 333373340
$ws.Range("I89").Value = 333373340
$ws.Range("K89").Value = 1666866700
$ws.Range("M89").Value = -1666861084

$ws.Range("H105").Value = 75003010
$ws.Range("I105").Value = 125000664
$ws.Range("J105").Value = 6521.375
$ws.Range("K105").Value = 125000664
$ws.Range("L105").Value = 6521.375
$ws.Range("M105").Value = -124998917
$ws.Range("N105").Value = -10015.375

$ws.Range("H107").Value = 7015.0527
$ws.Range("I107").Value = 7726.959
$ws.Range("K107").Value = 7726.959
$ws.Range("M107").Value = -5806.959

$ws.Range("H134").Value = 3980.3333
$ws.Range("J134").Value = 4333.3335
$ws.Range("L134").Value = 13000.0005
$ws.Range("N134").Value = -18070.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 937.90625
$ws.Range("I16").Value = 922.6087
$ws.Range("J16").Value = 977
$ws.Range("K16").Value = 922.6087
$ws.Range("L16").Value = 977
$ws.Range("M16").Value = -635.6087
$ws.Range("N16").Value = -1551

$ws.Range("H57").Value = 34999
$ws.Range("J57").Value = 34999
$ws.Range("L57").Value = 34999
$ws.Range("N57").Value = -36119

$ws.Range("H99").Value = 4427.364
$ws.Range("I99").Value = 3998.8
$ws.Range("K99").Value = 3998.8
$ws.Range("M99").Value = -2500.8

$ws.Range("H113").Value = 937.90625
$ws.Range("I113").Value = 922.6087
$ws.Range("J113").Value = 977
$ws.Range("K113").Value = 922.6087
$ws.Range("L113").Value = 977
$ws.Range("M113").Value = 1247.3913
$ws.Range("N113").Value = -5317

$ws.Range("H126").Value = 4427.364
$ws.Range("I126").Value = 3998.8
$ws.Range("K126").Value = 11996.4
$ws.Range("M126").Value = -9526.400000000001

$ws.Range("H132").Value = 2286.4167
$ws.Range("I132").Value = 1991.9474
$ws.Range("J132").Value = 3405.4
$ws.Range("K132").Value = 5975.8422
$ws.Range("L132").Value = 10216.2
$ws.Range("M132").Value = -3445.8422
$ws.Range("N132").Value = -15276.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J4").Value = 2237.5
$ws.Range("L4").Value = 6712.5
$ws.Range("N4").Value = -6936.5

$ws.Range("H14").Value = 2002.591
$ws.Range("I14").Value = 2002.591
$ws.Range("K14").Value = 6007.772999999999
$ws.Range("M14").Value = -5834.772999999999

$ws.Range("H22").Value = 850
$ws.Range("J22").Value = 1200
$ws.Range("L22").Value = 3600
$ws.Range("N22").Value = -3938

$ws.Range("H23").Value = 698.8333
$ws.Range("I23").Value = 665.6667
$ws.Range("J23").Value = 732
$ws.Range("K23").Value = 1997.0001
$ws.Range("L23").Value = 2196
$ws.Range("M23").Value = -1762.0001
$ws.Range("N23").Value = -2666

$ws.Range("H27").Value = 850
$ws.Range("J27").Value = 1200
$ws.Range("L27").Value = 3600
$ws.Range("N27").Value = -3804

$ws.Range("H100").Value = 5294.857
$ws.Range("I100").Value = 2666.5
$ws.Range("J100").Value = 7266.125
$ws.Range("K100").Value = 7999.5
$ws.Range("L100").Value = 21798.375
$ws.Range("M100").Value = -7188.5
$ws.Range("N100").Value = -23420.375

$ws.Range("H112").Value = 10237
$ws.Range("J112").Value = 14229.857
$ws.Range("L112").Value = 42689.571
$ws.Range("N112").Value = -44905.571

$ws.Range("H120").Value = 13016.5

$ws.Range("H121").Value = 8337350
$ws.Range("J121").Value = 16667758
$ws.Range("L121").Value = 50003274
$ws.Range("N121").Value = -50005894

$ws.Range("H127").Value = 45788.7
$ws.Range("I127").Value = 1500
$ws.Range("J127").Value = 50709.668
$ws.Range("K127").Value = 4500
$ws.Range("L127").Value = 152129.004
$ws.Range("M127").Value = 460
$ws.Range("N127").Value = -162049.004

$ws.Range("H131").Value = 1605
$ws.Range("J131").Value = 1628
$ws.Range("L131").Value = 4884
$ws.Range("N131").Value = -14964

$ws.Range("H134").Value = 4629.375
$ws.Range("I134").Value = 4360.25
$ws.Range("K134").Value = 13080.75
$ws.Range("M134").Value = -8010.75

$ws.Range("H137").Value = 1954.4048
$ws.Range("J137").Value = 1999.8649
$ws.Range("L137").Value = 5999.5947
$ws.Range("N137").Value = -16199.5947

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 44375.332
$ws.Range("J26").Value = 45042
$ws.Range("L26").Value = 45042
$ws.Range("N26").Value = -45602

$ws.Range("H50").Value = 44375.332
$ws.Range("J50").Value = 45042
$ws.Range("L50").Value = 45042
$ws.Range("N50").Value = -46038

$ws.Range("H113").Value = 30308772
$ws.Range("I113").Value = 40006716
$ws.Range("J113").Value = 2697.75
$ws.Range("K113").Value = 40006716
$ws.Range("L113").Value = 2697.75
$ws.Range("M113").Value = -40004546
$ws.Range("N113").Value = -7037.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11207.659
$ws.Range("I7").Value = 11599.568
$ws.Range("K7").Value = 11599.568
$ws.Range("M7").Value = -11487.568

$ws.Range("H16").Value = 1263.3125
$ws.Range("I16").Value = 1273.1936
$ws.Range("J16").Value = 957
$ws.Range("K16").Value = 1273.1936
$ws.Range("L16").Value = 957
$ws.Range("M16").Value = -1103.1936
$ws.Range("N16").Value = -1297

$ws.Range("H46").Value = 2504.0557
$ws.Range("J46").Value = 3419.5
$ws.Range("L46").Value = 3419.5
$ws.Range("N46").Value = -3795.5

$ws.Range("H55").Value = 1963.9445
$ws.Range("I55").Value = 506.1
$ws.Range("K55").Value = 506.1
$ws.Range("M55").Value = -333.1

$ws.Range("H126").Value = 11207.659
$ws.Range("I126").Value = 11599.568
$ws.Range("K126").Value = 34798.704
$ws.Range("M126").Value = -32328.704

$ws.Range("H132").Value = 2529999.2
$ws.Range("I132").Value = 3034591
$ws.Range("J132").Value = 7041.3335
$ws.Range("K132").Value = 9103773
$ws.Range("L132").Value = 21124.0005
$ws.Range("M132").Value = -9101243
$ws.Range("N132").Value = -26184.0005

$ws.Range("H136").Value = 7939421
$ws.Range("I136").Value = 11496634
$ws.Range("J136").Value = 4099.769
$ws.Range("K136").Value = 34489902
$ws.Range("L136").Value = 12299.307
$ws.Range("M136").Value = -34487352
$ws.Range("N136").Value = -17399.307

$ws.Range("H140").Value = 84768.60000000001
$ws.Range("J140").Value = 84768.60000000001
$ws.Range("L140").Value = 84768.60000000001
$ws.Range("N140").Value = -95128.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 63374.6
$ws.Range("J123").Value = 63374.6
$ws.Range("L123").Value = 63374.6
$ws.Range("N123").Value = -73174.60000000001

$ws.Range("H135").Value = 68428.14
$ws.Range("J135").Value = 68428.14
$ws.Range("L135").Value = 68428.14
$ws.Range("N135").Value = -78568.14

$ws.Range("H136").Value = 3498097.5
$ws.Range("I136").Value = 4275072
$ws.Range("J136").Value = 1712
$ws.Range("K136").Value = 12825216
$ws.Range("L136").Value = 5136
$ws.Range("M136").Value = -12822666
$ws.Range("N136").Value = -10236

Write-Host "Applied all cell updates"